$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header / data cells with new shared-string values
$ws.Range("A1").Value = "IMPORT NHÂN SỰ"
$ws.Range("C2").Value = "Mã nhân sự"
$ws.Range("D2").Value = "Tên phòng"

$ws.Range("B4").Value = "NS71"
$ws.Range("C4").Value = "NS1"
$ws.Range("D4").Value = "Phòng tài chính"

$ws.Range("B5").Value = "NS42"
$ws.Range("C5").Value = "NS5"
$ws.Range("D5").Value = "Phòng kinh doanh"

$ws.Range("B6").Value = "NS15"
$ws.Range("C6").Value = "NS6"
$ws.Range("D6").Value = "Phòng kinh doanh"

$ws.Range("B7").Value = "NS61"
$ws.Range("C7").Value = "NS7"
$ws.Range("D7").Value = "Phòng nhân sự"

# Move active selection to C9
$ws.Range("C9").Select()
